$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the 2021 column (R) to the maternal-mortality table, mirroring the
# formatting already used for 2020 (column Q).
$ws.Range("Q3:Q14").Copy()
$ws.Range("R3:R14").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Header year
$ws.Range("R4").Value = 2021

# Maternal mortality rate per 100,000 live births, by region, for 2021
$ws.Range("R5").Value = 33.3
$ws.Range("R6").Value = 38.3
$ws.Range("R7").Value = 31.7
$ws.Range("R8").Value = 98.7
$ws.Range("R9").Value = 157.2
$ws.Range("R10").Value = 24.9
$ws.Range("R11").Value = 38.4
$ws.Range("R12").Value = 15.1
$ws.Range("R13").Value = 14.6
$ws.Range("R14").Value = 21.7

# Match the author's final cursor position / selection
$ws.Range("S6").Select() | Out-Null
